# Insert 3 new data rows (weekly update) right above the existing row 569,
# which pushes the existing rows 569..673 down to 572..676.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("569:571").Insert()

# Constant values shared by every data row in this sheet
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$catId     = 100112023
$categoria = "Brócoli"
$variedad  = "Sin especificar"
$unidad    = "`$/unidad"
$origen    = "Región de Arica y Parinacota"
$kgUnid    = 1
$clasif    = "Hortaliza"

function Set-DataRow($RowIndex, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm) {
    $ws.Cells.Item($RowIndex, 1).Value  = $mercadoId
    $ws.Cells.Item($RowIndex, 2).Value  = $mercado
    $ws.Cells.Item($RowIndex, 3).Value  = $region
    $ws.Cells.Item($RowIndex, 4).Value  = $Fecha
    $ws.Cells.Item($RowIndex, 5).Value  = $codreg
    $ws.Cells.Item($RowIndex, 6).Value  = $catId
    $ws.Cells.Item($RowIndex, 7).Value  = $categoria
    $ws.Cells.Item($RowIndex, 8).Value  = $variedad
    $ws.Cells.Item($RowIndex, 9).Value  = $Calidad
    $ws.Cells.Item($RowIndex, 10).Value = $Volumen
    $ws.Cells.Item($RowIndex, 11).Value = $PrecioMin
    $ws.Cells.Item($RowIndex, 12).Value = $PrecioMax
    $ws.Cells.Item($RowIndex, 13).Value = $PrecioProm
    $ws.Cells.Item($RowIndex, 14).Value = $unidad
    $ws.Cells.Item($RowIndex, 15).Value = $origen
    $ws.Cells.Item($RowIndex, 16).Value = $PrecioProm
    $ws.Cells.Item($RowIndex, 17).Value = $kgUnid
    $ws.Cells.Item($RowIndex, 18).Value = $clasif
}

Set-DataRow 569 45244 "Primera" 1000 700 800 750
Set-DataRow 570 45244 "Segunda" 1200 500 600 550
Set-DataRow 571 45244 "Tercera" 1200 350 400 375
